# Update scripts with new TPM-derived NATMI ligand-receptor numbers.
# The underlying Python pipeline re-ran with refreshed TPM expression
# values, so the cached ligand/receptor/edge expression + specificity
# metrics (columns G-J, M-T) change for every Sending/Target cluster
# combination on the sheet. Columns A-F, K and L (identifiers, counts
# and detection rates) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.922976999999999
$ws.Range("H2").Value = 14.768931
$ws.Range("I2").Value = 0.2897120038548413
$ws.Range("J2").Value = 0.2897120038548412
$ws.Range("M2").Value = 0.3822983333333334
$ws.Range("N2").Value = 1.146895
$ws.Range("O2").Value = 0.04915201160487953
$ws.Range("P2").Value = 0.04915201160487953
$ws.Range("Q2").Value = 1.882045902138333
$ws.Range("R2").Value = 16.938413119245
$ws.Range("S2").Value = 0.01423992777554606
$ws.Range("T2").Value = 0.01423992777554606
$ws.Range("G3").Value = 4.922976999999999
$ws.Range("H3").Value = 14.768931
$ws.Range("I3").Value = 0.2897120038548413
$ws.Range("J3").Value = 0.2897120038548412
$ws.Range("O3").Value = 0.3087451919724631
$ws.Range("P3").Value = 0.3087451919724631
$ws.Range("Q3").Value = 11.82194999520633
$ws.Range("R3").Value = 106.397549956857
$ws.Range("S3").Value = 0.08944718824688996
$ws.Range("T3").Value = 0.08944718824688994
$ws.Range("G4").Value = 4.922976999999999
$ws.Range("H4").Value = 14.768931
$ws.Range("I4").Value = 0.2897120038548413
$ws.Range("J4").Value = 0.2897120038548412
$ws.Range("O4").Value = 0.6421027964226573
$ws.Range("P4").Value = 0.6421027964226573
$ws.Range("Q4").Value = 24.586316964469
$ws.Range("R4").Value = 221.276852680221
$ws.Range("S4").Value = 0.1860248878324053
$ws.Range("T4").Value = 0.1860248878324053
$ws.Range("I5").Value = 0.5769489387710858
$ws.Range("J5").Value = 0.5769489387710858
$ws.Range("M5").Value = 0.3822983333333334
$ws.Range("N5").Value = 1.146895
$ws.Range("O5").Value = 0.04915201160487953
$ws.Range("P5").Value = 0.04915201160487953
$ws.Range("Q5").Value = 3.748013100973334
$ws.Range("R5").Value = 33.73211790876
$ws.Range("S5").Value = 0.02835820093389934
$ws.Range("T5").Value = 0.02835820093389934
$ws.Range("I6").Value = 0.5769489387710858
$ws.Range("J6").Value = 0.5769489387710858
$ws.Range("O6").Value = 0.3087451919724631
$ws.Range("P6").Value = 0.3087451919724631
$ws.Range("S6").Value = 0.1781302108591878
$ws.Range("T6").Value = 0.1781302108591878
$ws.Range("I7").Value = 0.5769489387710858
$ws.Range("J7").Value = 0.5769489387710858
$ws.Range("O7").Value = 0.6421027964226573
$ws.Range("P7").Value = 0.6421027964226573
$ws.Range("S7").Value = 0.3704605269779987
$ws.Range("T7").Value = 0.3704605269779987
$ws.Range("I8").Value = 0.133339057374073
$ws.Range("J8").Value = 0.133339057374073
$ws.Range("M8").Value = 0.3822983333333334
$ws.Range("N8").Value = 1.146895
$ws.Range("O8").Value = 0.04915201160487953
$ws.Range("P8").Value = 0.04915201160487953
$ws.Range("Q8").Value = 0.8662058291916668
$ws.Range("R8").Value = 7.795852462725001
$ws.Range("S8").Value = 0.006553882895434133
$ws.Range("T8").Value = 0.006553882895434133
$ws.Range("I9").Value = 0.133339057374073
$ws.Range("J9").Value = 0.133339057374073
$ws.Range("O9").Value = 0.3087451919724631
$ws.Range("P9").Value = 0.3087451919724631
$ws.Range("Q9").Value = 5.441016070131668
$ws.Range("R9").Value = 48.96914463118501
$ws.Range("S9").Value = 0.04116779286638544
$ws.Range("T9").Value = 0.04116779286638544
$ws.Range("I10").Value = 0.133339057374073
$ws.Range("J10").Value = 0.133339057374073
$ws.Range("O10").Value = 0.6421027964226573
$ws.Range("P10").Value = 0.6421027964226573
$ws.Range("S10").Value = 0.08561738161225341
$ws.Range("T10").Value = 0.08561738161225341
